$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.068.23"
Set-TextValue $ws.Range("E2") "  -1.15%  "
Set-TextValue $ws.Range("D3") "1.664.81"
Set-TextValue $ws.Range("E3") "  -1.35%  "
Set-TextValue $ws.Range("D5") "209.47"
Set-TextValue $ws.Range("E5") "  -4.16%  "
Set-TextValue $ws.Range("D6") "0.5159"
Set-TextValue $ws.Range("E6") "  -5.72%  "
Set-TextValue $ws.Range("D7") "1.004"
Set-TextValue $ws.Range("E7") "  -0.70%  "
Set-TextValue $ws.Range("D8") "0.2635"
Set-TextValue $ws.Range("E8") "  -3.13%  "
Set-TextValue $ws.Range("D9") "0.06209"
Set-TextValue $ws.Range("E9") "  -3.87%  "
Set-TextValue $ws.Range("D10") "20.85"
Set-TextValue $ws.Range("E10") "  -5.33%  "
Set-TextValue $ws.Range("D11") "0.07504"
Set-TextValue $ws.Range("E11") "  -2.36%  "
Set-TextValue $ws.Range("D12") "1.671.09"
Set-TextValue $ws.Range("E12") "  -1.38%  "
Set-TextValue $ws.Range("D13") "4.415"
Set-TextValue $ws.Range("E13") "  -2.54%  "
Set-TextValue $ws.Range("D14") "0.5568"
Set-TextValue $ws.Range("E14") "  -4.22%  "
Set-TextValue $ws.Range("D15") "0.000007948"
Set-TextValue $ws.Range("E15") "  -5.06%  "
Set-TextValue $ws.Range("D16") "65.30"
Set-TextValue $ws.Range("E16") "  +0.37%  "
Set-TextValue $ws.Range("D17") "26.092.93"
Set-TextValue $ws.Range("E17") "  -1.22%  "
Set-TextValue $ws.Range("D18") "1.004"
Set-TextValue $ws.Range("E18") "  -0.63%  "
Set-TextValue $ws.Range("D19") "4.775"
Set-TextValue $ws.Range("E19") "  -3.29%  "
Set-TextValue $ws.Range("D20") "10.36"
Set-TextValue $ws.Range("E20") "  -5.59%  "
Set-TextValue $ws.Range("D21") "185.57"
Set-TextValue $ws.Range("E21") "  -2.58%  "
Set-TextValue $ws.Range("D22") "6.131"
Set-TextValue $ws.Range("E22") "  -1.63%  "
Set-TextValue $ws.Range("D23") "1.004"
Set-TextValue $ws.Range("E23") "  -0.74%  "
Set-TextValue $ws.Range("D24") "146.12"
Set-TextValue $ws.Range("E24") "  -2.34%  "
Set-TextValue $ws.Range("D25") "0.1243"
Set-TextValue $ws.Range("E25") "  -5.73%  "
Set-TextValue $ws.Range("D26") "7.535"
Set-TextValue $ws.Range("E26") "  -4.27%  "
Set-TextValue $ws.Range("D27") "15.67"
Set-TextValue $ws.Range("E27") "  -0.18%  "
Set-TextValue $ws.Range("D28") "0.06330"
Set-TextValue $ws.Range("E28") "  +0.15%  "
Set-TextValue $ws.Range("D29") "1.338"
Set-TextValue $ws.Range("E29") "  -5.06%  "
Set-TextValue $ws.Range("D30") "1.271"
Set-TextValue $ws.Range("E30") "  -4.33%  "
Set-TextValue $ws.Range("D31") "3.468"
Set-TextValue $ws.Range("E31") "  -3.04%  "
Set-TextValue $ws.Range("D32") "3.421"
Set-TextValue $ws.Range("E32") "  -4.48%  "
Set-TextValue $ws.Range("D33") "1.614"
Set-TextValue $ws.Range("E33") "  -3.64%  "
Set-TextValue $ws.Range("D34") "0.9926"
Set-TextValue $ws.Range("E34") "  -4.58%  "
Set-TextValue $ws.Range("D35") "2.408"
Set-TextValue $ws.Range("E35") "  -0.11%  "
Set-TextValue $ws.Range("D36") "0.6019"
Set-TextValue $ws.Range("E36") "  -2.23%  "
Set-TextValue $ws.Range("D37") "2.701"
Set-TextValue $ws.Range("E37") "  -0.59%  "
Set-TextValue $ws.Range("D38") "6.069"
Set-TextValue $ws.Range("E38") "  -2.72%  "
Set-TextValue $ws.Range("D39") "0.01603"
Set-TextValue $ws.Range("E39") "  -1.26%  "
Set-TextValue $ws.Range("D40") "1.077.87"
Set-TextValue $ws.Range("E40") "  -3.14%  "
Set-TextValue $ws.Range("D41") "0.8588"
Set-TextValue $ws.Range("E41") "  -2.57%  "
Set-TextValue $ws.Range("E42") "  -1.19%  "
Set-TextValue $ws.Range("D43") "99.04"
Set-TextValue $ws.Range("E43") "  -2.32%  "
Set-TextValue $ws.Range("D44") "1.812.85"
Set-TextValue $ws.Range("E44") "  -1.43%  "
Set-TextValue $ws.Range("D45") "0.00000000111"
Set-TextValue $ws.Range("E45") "  +0.91%  "
Set-TextValue $ws.Range("D46") "55.99"
Set-TextValue $ws.Range("E46") "  -2.31%  "
Set-TextValue $ws.Range("E47") "  -0.35%  "
Set-TextValue $ws.Range("D48") "0.05250"
Set-TextValue $ws.Range("D49") "7.923"
Set-TextValue $ws.Range("E49") "  -3.07%  "
Set-TextValue $ws.Range("D50") "0.4259"
Set-TextValue $ws.Range("E50") "  -1.04%  "
Set-TextValue $ws.Range("D51") "5.886"
Set-TextValue $ws.Range("E51") "  -2.52%  "
